# Auto-applied updates to Leve profit calculation sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 446.54544
$ws.Range("I55").Value = 544.2857
$ws.Range("J55").Value = 275.5
$ws.Range("K55").Value = 544.2857
$ws.Range("L55").Value = 275.5
$ws.Range("M55").Value = -330.2857
$ws.Range("N55").Value = -703.5
$ws.Range("H64").Value = 5550
$ws.Range("I64").Value = 7000
$ws.Range("J64").Value = 3520
$ws.Range("K64").Value = 7000
$ws.Range("L64").Value = 3520
$ws.Range("M64").Value = -6752
$ws.Range("N64").Value = -4016
$ws.Range("H67").Value = 5550
$ws.Range("I67").Value = 7000
$ws.Range("J67").Value = 3520
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 3520
$ws.Range("M67").Value = -6142
$ws.Range("N67").Value = -5236
$ws.Range("H76").Value = 4634625
$ws.Range("I76").Value = 11115200
$ws.Range("K76").Value = 11115200
$ws.Range("M76").Value = -11114885
$ws.Range("H79").Value = 4634625
$ws.Range("I79").Value = 11115200
$ws.Range("K79").Value = 11115200
$ws.Range("M79").Value = -11114108
$ws.Range("H138").Value = 3403.7527
$ws.Range("I138").Value = 1705.3871
$ws.Range("J138").Value = 4252.9355
$ws.Range("K138").Value = 5116.1613
$ws.Range("L138").Value = 12758.8065
$ws.Range("M138").Value = 23.83870000000024
$ws.Range("N138").Value = -23038.8065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 214697.12
$ws.Range("I61").Value = 1733.4872
$ws.Range("J61").Value = 1252894.9
$ws.Range("K61").Value = 1733.4872
$ws.Range("L61").Value = 1252894.9
$ws.Range("M61").Value = -1521.4872
$ws.Range("N61").Value = -1253318.9
$ws.Range("H135").Value = 62645.445
$ws.Range("J135").Value = 62645.445
$ws.Range("L135").Value = 62645.445
$ws.Range("N135").Value = -72785.44500000001
$ws.Range("H136").Value = 214697.12
$ws.Range("I136").Value = 1733.4872
$ws.Range("J136").Value = 1252894.9
$ws.Range("K136").Value = 5200.461600000001
$ws.Range("L136").Value = 3758684.7
$ws.Range("M136").Value = -2650.461600000001
$ws.Range("N136").Value = -3763784.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 63268.668
$ws.Range("J135").Value = 63268.668
$ws.Range("L135").Value = 63268.668
$ws.Range("N135").Value = -73408.66800000001
$ws.Range("H140").Value = 44570.91
$ws.Range("J140").Value = 47028
$ws.Range("L140").Value = 47028
$ws.Range("N140").Value = -57388

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 945.375
$ws.Range("I5").Value = 687
$ws.Range("J5").Value = 1376
$ws.Range("K5").Value = 687
$ws.Range("L5").Value = 1376
$ws.Range("M5").Value = -575
$ws.Range("N5").Value = -1600
$ws.Range("H8").Value = 2567.5
$ws.Range("I8").Value = 90
$ws.Range("J8").Value = 3393.3333
$ws.Range("K8").Value = 90
$ws.Range("L8").Value = 3393.3333
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = -3673.3333
$ws.Range("H10").Value = 50008
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H11").Value = 950
$ws.Range("I11").Value = 1100
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 1100
$ws.Range("L11").Value = 800
$ws.Range("M11").Value = -960
$ws.Range("N11").Value = -1080
$ws.Range("H13").Value = 26729.092
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 28902
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 28902
$ws.Range("M13").Value = -4861
$ws.Range("N13").Value = -29180
$ws.Range("H15").Value = 360
$ws.Range("I15").Value = 360
$ws.Range("K15").Value = 360
$ws.Range("M15").Value = -190
$ws.Range("H19").Value = 666.3333
$ws.Range("I19").Value = 666.3333
$ws.Range("K19").Value = 666.3333
$ws.Range("M19").Value = -496.3333
$ws.Range("H24").Value = 666.3333
$ws.Range("I24").Value = 666.3333
$ws.Range("K24").Value = 666.3333
$ws.Range("M24").Value = -496.3333
$ws.Range("H26").Value = 7220
$ws.Range("J26").Value = 7220
$ws.Range("L26").Value = 7220
$ws.Range("N26").Value = -7794
$ws.Range("H31").Value = 6498294
$ws.Range("I31").Value = 1437.6571
$ws.Range("J31").Value = 11912341
$ws.Range("K31").Value = 1437.6571
$ws.Range("L31").Value = 11912341
$ws.Range("M31").Value = -1142.6571
$ws.Range("N31").Value = -11912931
$ws.Range("H34").Value = 6498294
$ws.Range("I34").Value = 1437.6571
$ws.Range("J34").Value = 11912341
$ws.Range("K34").Value = 1437.6571
$ws.Range("L34").Value = 11912341
$ws.Range("M34").Value = -1235.6571
$ws.Range("N34").Value = -11912745
$ws.Range("H122").Value = 4037.875
$ws.Range("I122").Value = 3648.6924
$ws.Range("J122").Value = 4497.8184
$ws.Range("K122").Value = 10946.0772
$ws.Range("L122").Value = 13493.4552
$ws.Range("M122").Value = -8496.0772
$ws.Range("N122").Value = -18393.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 172906.83
$ws.Range("I113").Value = 482.5909
$ws.Range("K113").Value = 1447.7727
$ws.Range("M113").Value = 722.2273
$ws.Range("H131").Value = 2084359.9
$ws.Range("I131").Value = 6250684.5
$ws.Range("J131").Value = 1197.625
$ws.Range("K131").Value = 18752053.5
$ws.Range("L131").Value = 3592.875
$ws.Range("M131").Value = -18747013.5
$ws.Range("N131").Value = -13672.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9586.429
$ws.Range("I80").Value = 10725.833
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 10725.833
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -9727.833000000001
$ws.Range("N80").Value = -4746
$ws.Range("H83").Value = 9586.429
$ws.Range("I83").Value = 10725.833
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 53629.165
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -48637.165
$ws.Range("N83").Value = -23734
$ws.Range("H140").Value = 34539.57
$ws.Range("J140").Value = 34539.57
$ws.Range("L140").Value = 34539.57
$ws.Range("N140").Value = -44899.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2822.2222
$ws.Range("I40").Value = 2800
$ws.Range("K40").Value = 2800
$ws.Range("M40").Value = -2664
$ws.Range("H61").Value = 1339.8
$ws.Range("I61").Value = 1156
$ws.Range("J61").Value = 2075
$ws.Range("K61").Value = 1156
$ws.Range("L61").Value = 2075
$ws.Range("M61").Value = -954
$ws.Range("N61").Value = -2479
$ws.Range("H113").Value = 1339.8
$ws.Range("I113").Value = 1156
$ws.Range("J113").Value = 2075
$ws.Range("K113").Value = 1156
$ws.Range("L113").Value = 2075
$ws.Range("M113").Value = 1014
$ws.Range("N113").Value = -6415
$ws.Range("H134").Value = 59329.77
$ws.Range("J134").Value = 59329.77
$ws.Range("L134").Value = 59329.77
$ws.Range("N134").Value = -69469.76999999999
$ws.Range("H140").Value = 107214.5
$ws.Range("J140").Value = 107214.5
$ws.Range("L140").Value = 107214.5
$ws.Range("N140").Value = -117574.5
$ws.Range("H141").Value = 92070.7
$ws.Range("J141").Value = 77856.336
$ws.Range("L141").Value = 77856.336
$ws.Range("N141").Value = -88216.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 50000784
$ws.Range("I107").Value = 71429260
$ws.Range("K107").Value = 214287780
$ws.Range("M107").Value = -214285860
$ws.Range("H132").Value = 1714.766
$ws.Range("I132").Value = 867.2593000000001
$ws.Range("J132").Value = 2858.9
$ws.Range("K132").Value = 2601.7779
$ws.Range("L132").Value = 8576.700000000001
$ws.Range("M132").Value = -71.77790000000005
$ws.Range("N132").Value = -13636.7
